$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- Column A: Time (rerun timestamps) ---
$ws.Range("A2").Value = "20160403_214048"
$ws.Range("A3").Value = "20160403_214226"
$ws.Range("A4").Value = "20160403_214442"
$ws.Range("A5").Value = "20160403_214726"
$ws.Range("A6").Value = "20160403_214955"
$ws.Range("A7").Value = "20160403_220626"
$ws.Range("A8").Value = "20160403_225448"
$ws.Range("A9").Value = "20160403_230139"
$ws.Range("A10").Value = "20160403_230823"
$ws.Range("A11").Value = "20160403_231456"

# --- Column B: RunningTime(s) ---
$ws.Range("B2").Value = 98.201
$ws.Range("B3").Value = 135.334
$ws.Range("B4").Value = 164.233
$ws.Range("B5").Value = 149.209
$ws.Range("B6").Value = 125.242
$ws.Range("B7").Value = 2901.521
$ws.Range("B8").Value = 410.816
$ws.Range("B9").Value = 403.969
$ws.Range("B10").Value = 393.342
$ws.Range("B11").Value = 411.602

# --- Column C: Preprocess description (reordered steps) ---
$ws.Range("C2:C6").Value = 'remove multiple spaces, trim "space" and ",", convert to lower, convert unicode to ascii'
$ws.Range("C7:C11").Value = 'trim "space" and ",", convert unicode to ascii, convert to lower, remove multiple spaces'

# --- Column G: Test_Accuracy ---
$ws.Range("G2").Value = 0.956
$ws.Range("G3").Value = 0.967333333333333
$ws.Range("G4").Value = 0.964
$ws.Range("G5").Value = 0.959333333333333
$ws.Range("G7").Value = 0.957333333333333
$ws.Range("G8").Value = 0.967333333333333
$ws.Range("G9").Value = 0.955333333333333
$ws.Range("G10").Value = 0.955333333333333
$ws.Range("G11").Value = 0.956666666666667

# --- Column J ---
$ws.Range("J2").Value = 0.163265306122449
$ws.Range("J3").Value = 0.163265306122449
$ws.Range("J4").Value = 0.163265306122449
$ws.Range("J5").Value = 0.153061224489796
$ws.Range("J6").Value = 0.173469387755102
$ws.Range("J7").Value = 0.173469387755102
$ws.Range("J8").Value = 0.173469387755102
$ws.Range("J9").Value = 0.173469387755102
$ws.Range("J10").Value = 0.163265306122449
$ws.Range("J11").Value = 0.163265306122449
